# Weekly update: Fruta / hortaliza, semanal
# Two new weekly price records are inserted at rows 32-33 (pushing all
# historical rows previously at 32-117 down to 34-119), and the new
# records are populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 32, shifting the
# existing data (previously rows 32-117) down to rows 34-119.
$ws.Range("A32:T33").EntireRow.Insert()

# Populate the first new record (row 32).
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = 'Femacal de La Calera'
$ws.Range("C32").Value = 'Coquimbo'
$ws.Range("D32").Value = 44487
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 'Fruta'
$ws.Range("G32").Value = 100101
$ws.Range("H32").Value = 'Berries'
$ws.Range("I32").Value = 100101001
$ws.Range("J32").Value = 'Arándano (blue)'
$ws.Range("K32").Value = 'Sin especificar'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 128
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 11000
$ws.Range("P32").Value = 10531
$ws.Range("Q32").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R32").Value = 'Provincia de Quillota'
$ws.Range("S32").Value = 7021
$ws.Range("T32").Value = 1.5

# Populate the second new record (row 33).
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = 'Femacal de La Calera'
$ws.Range("C33").Value = 'Coquimbo'
$ws.Range("D33").Value = 44487
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 'Fruta'
$ws.Range("G33").Value = 100101
$ws.Range("H33").Value = 'Berries'
$ws.Range("I33").Value = 100101001
$ws.Range("J33").Value = 'Arándano (blue)'
$ws.Range("K33").Value = 'Sin especificar'
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 48
$ws.Range("N33").Value = 25000
$ws.Range("O33").Value = 25000
$ws.Range("P33").Value = 25000
$ws.Range("Q33").Value = '$/bandeja 5 kilos'
$ws.Range("R33").Value = 'Provincia de Quillota'
$ws.Range("S33").Value = 5000
$ws.Range("T33").Value = 5
